# updated print to excel functions and ship parallels
#
# 1. Rename Sheet1 -> system_structures, add a "complex sys 2" row to it.
# 2. Add a new sheet "ship_structure" (right after system_structures) that
#    mirrors the "complex sys" print-to-excel pattern used elsewhere in the
#    workbook, holding the ship's own structure data.
# 3. Leave final selection/zoom/active-sheet state matching the author's
#    last interactive session.

$wb = $excel.ActiveWorkbook

# --- rename the second sheet ---------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "system_structures"

# --- new "ship_structure" sheet, inserted right after system_structures --
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws2)
$ws3.Name = "ship_structure"
$ws3.Range("A1").Value = "ship structure"
$ws3.Range("A2").Value = "[(1,2), (3,4)]"

# --- append the new "complex sys 2" row to system_structures -------------
$ws2.Range("A5").Value = "complex sys 2"
$ws2.Range("B5").Value = "[([0,4], [1,5]), (6,7)]"

# --- restore view/selection state -----------------------------------------
$ws3.Range("A3").Select() | Out-Null

$ws2.Activate() | Out-Null
$ws2.Range("C9").Select() | Out-Null
$excel.ActiveWindow.Zoom = 235
